# The "keyboard row" leetcode problem (#500) is added as a new row to the
# "哈希" (hash) sheet, which is the workbook's active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (existing last data row) already carries the table's formatting
# (style index 4, wrap text, etc). Clone that formatting down into the new
# row 10 (new data row) and row 11 (blank spacer row that follows it, just
# like the existing sheet layout) before writing any values.
$ws.Range("A9:H9").Copy()
$ws.Range("A10:H11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 10 content: problem #9 / leetcode 500 ("Keyboard Row").
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 500

# Populate the new unique shared strings in the same order they are first
# used left-to-right in the row (time/space complexity columns first, then
# the problem statement, then the solution write-up) so the generated
# shared-string table lines up with the target ordering.
$ws.Range("F10").Value = "O(N*M), N是数组长度，M是字符串长度"
$ws.Range("G10").Value = "O(N*M), N是数组长度，M是字符串长度"
$ws.Range("C10").Value = "给定一个单词列表，只返回可以使用在键盘同一行的字母打印出来的单词。键盘如下图所示。 "
$ws.Range("D10").Value = "1 获取键盘中三行字符的位置，自己构建map`n2 解析字符串的字符，保留其位置，迭代下一个字符，并与prev比较，`n3 如果不相等就说明不相登，false`n4 如果相等就继续比较"
$ws.Range("E10").Value = "哈希表`n字符出现次数"

# Row heights: the new data row is tall (110pt, matching other multi-line
# rows) and row 11 is a short blank spacer row (21pt) under it.
$ws.Rows.Item(10).RowHeight = 110
$ws.Rows.Item(11).RowHeight = 21

# Match the author's final cursor position/selection in the sheet.
$ws.Range("D12").Select()
